$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Panel A (Bond Futures), FF1, Avg Daily Volume ---
$ws.Range("E2").Value = 6953.776198605728

# --- Row 26: Panel B (E-mini Futures), Emini, Avg Daily Volume ---
$ws.Range("D26").Value = 898789.6108597285
$ws.Range("E26").Value = 864606.3589447475
$ws.Range("F26").Value = 2468
$ws.Range("G26").Value = 1002885
$ws.Range("H26").Value = 1406731
$ws.Range("I26").Value = 221
$ws.Range("J26").Value = 1157603.194570136
$ws.Range("K26").Value = 956073.8115222843
$ws.Range("L26").Value = 489073
$ws.Range("M26").Value = 1178049
$ws.Range("N26").Value = 1614520
$ws.Range("O26").Value = 221
$ws.Range("V26").Value = 1400512.303167421
$ws.Range("W26").Value = 1053782.529173701
$ws.Range("X26").Value = 568235
$ws.Range("Y26").Value = 1553653
$ws.Range("Z26").Value = 1984535
$ws.Range("AA26").Value = 221
$ws.Range("AB26").Value = 1267942.162895928
$ws.Range("AC26").Value = 1034806.083474207
$ws.Range("AD26").Value = 0
$ws.Range("AE26").Value = 1324144
$ws.Range("AF26").Value = 1916380
$ws.Range("AG26").Value = 221

# --- Row 27: Panel B (E-mini Futures), Emini, Diff_Vol (Ann - Day) ---
$ws.Range("D27").Value = 430664.479638009
$ws.Range("J27").Value = 171850.8959276018
$ws.Range("V27").Value = -71058.21266968326
$ws.Range("AB27").Value = 61511.92760180996

# --- Row 28: Panel B (E-mini Futures), Emini, # Obs ---
$ws.Range("D28").Value = 221
$ws.Range("J28").Value = 221
$ws.Range("V28").Value = 221
$ws.Range("AB28").Value = 221
